$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Funktionen")
$ws.Rows.Item(7).Delete()
$ws.Range("A7").Select()
